$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for the 3 new rows by copying from the last existing row
$ws.Range("A16:M16").Copy($ws.Range("A17:M17"))
$ws.Range("A16:M16").Copy($ws.Range("A18:M18"))
$ws.Range("A16:M16").Copy($ws.Range("A19:M19"))

# Overwrite cell values for rows 10 through 19 with the new averaged-intensity results
# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.674365088669574
$ws.Range("D10").Value = 1.549275029639975
$ws.Range("E10").Value = 0.6158939211221109
$ws.Range("F10").Value = 1.674365088669574
$ws.Range("G10").Value = 1.228437297018139
$ws.Range("H10").Value = 1.20694505993909
$ws.Range("I10").Value = 0.7380088131999583
$ws.Range("J10").Value = 1.549275029639975
$ws.Range("K10").Value = 1.082584475381043
$ws.Range("L10").Value = 1.378474782025308
$ws.Range("M10").Value = 1.168820868264808

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.8641245574885912
$ws.Range("D11").Value = 1.770979305668895
$ws.Range("E11").Value = 1.113894330462043
$ws.Range("F11").Value = 0.8641245574885912
$ws.Range("G11").Value = 0.5825509544718723
$ws.Range("H11").Value = 2.729480017521278
$ws.Range("I11").Value = 0.861675581484458
$ws.Range("J11").Value = 1.770979305668895
$ws.Range("K11").Value = 1.442436818065469
$ws.Range("L11").Value = 1.15328068777703
$ws.Range("M11").Value = 1.320450791182856

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.8653353956300618
$ws.Range("D12").Value = 1.779038127328181
$ws.Range("E12").Value = 1.113416728255973
$ws.Range("F12").Value = 0.8653353956300618
$ws.Range("G12").Value = 0.5848169410114502
$ws.Range("H12").Value = 2.723699177121625
$ws.Range("I12").Value = 0.8607404800565623
$ws.Range("J12").Value = 1.779038127328181
$ws.Range("K12").Value = 1.446227427792077
$ws.Range("L12").Value = 1.155781411711069
$ws.Range("M12").Value = 1.321174474900642

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.8634874956247396
$ws.Range("D13").Value = 1.775932319162703
$ws.Range("E13").Value = 1.113567780979181
$ws.Range("F13").Value = 0.8634874956247396
$ws.Range("G13").Value = 0.5837122533722552
$ws.Range("H13").Value = 2.729666530107401
$ws.Range("I13").Value = 0.8611096382837026
$ws.Range("J13").Value = 1.775932319162703
$ws.Range("K13").Value = 1.444750050070942
$ws.Range("L13").Value = 1.154118772847841
$ws.Range("M13").Value = 1.321246002921664

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.1746239999999965
$ws.Range("D14").Value = 1.371964000000002
$ws.Range("E14").Value = 1.530551999999999
$ws.Range("F14").Value = 0.1746239999999965
$ws.Range("G14").Value = 0.8761240000000011
$ws.Range("H14").Value = 1.011708000000002
$ws.Range("I14").Value = 0.8871640000000003
$ws.Range("J14").Value = 1.371964000000002
$ws.Range("K14").Value = 1.451258000000001
$ws.Range("L14").Value = 0.8129409999999985
$ws.Range("M14").Value = 0.9753560000000002

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.01
$ws.Range("D15").Value = 0.17
$ws.Range("E15").Value = 1.886687500000001
$ws.Range("F15").Value = 0.01
$ws.Range("G15").Value = 0.8573874999999987
$ws.Range("H15").Value = 0.03
$ws.Range("I15").Value = 0.9938999999999982
$ws.Range("J15").Value = 0.17
$ws.Range("K15").Value = 1.028343750000001
$ws.Range("L15").Value = 0.5191718750000003
$ws.Range("M15").Value = 0.657995833333333

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.3993886746624021
$ws.Range("D16").Value = 0.5165550780416036
$ws.Range("E16").Value = 1.492762168012801
$ws.Range("F16").Value = 0.3993886746624021
$ws.Range("G16").Value = 0.9448966324224
$ws.Range("H16").Value = 0.4362784116736026
$ws.Range("I16").Value = 0.9905982595072008
$ws.Range("J16").Value = 0.5165550780416036
$ws.Range("K16").Value = 1.004658623027202
$ws.Range("L16").Value = 0.7020236488448023
$ws.Range("M16").Value = 0.7967465373866683

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9875976918588665
$ws.Range("D17").Value = 0.9956729380976406
$ws.Range("E17").Value = 0.9929618368150657
$ws.Range("F17").Value = 0.9875976918588665
$ws.Range("G17").Value = 0.9972591882637242
$ws.Range("H17").Value = 0.9916696337485124
$ws.Range("I17").Value = 0.988845356407998
$ws.Range("J17").Value = 0.9956729380976406
$ws.Range("K17").Value = 0.9943173874563531
$ws.Range("L17").Value = 0.9909575396576098
$ws.Range("M17").Value = 0.9923344408653012

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.194908510254233
$ws.Range("D18").Value = 0.954400535753294
$ws.Range("E18").Value = 1.020073377622325
$ws.Range("F18").Value = 1.194908510254233
$ws.Range("G18").Value = 0.9562141094853609
$ws.Range("H18").Value = 1.002062456044293
$ws.Range("I18").Value = 0.9524729883798714
$ws.Range("J18").Value = 0.954400535753294
$ws.Range("K18").Value = 0.9872369566878096
$ws.Range("L18").Value = 1.091072733471021
$ws.Range("M18").Value = 1.013355329589896

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9333259351345736
$ws.Range("D19").Value = 0.7558880156610238
$ws.Range("E19").Value = 1.050042553223777
$ws.Range("F19").Value = 0.9333259351345736
$ws.Range("G19").Value = 0.9538852878385292
$ws.Range("H19").Value = 1.029393550988513
$ws.Range("I19").Value = 1.023314237701116
$ws.Range("J19").Value = 0.7558880156610238
$ws.Range("K19").Value = 0.9029652844424005
$ws.Range("L19").Value = 0.9181456097884871
$ws.Range("M19").Value = 0.9576415967579223
